# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (GitHub Actions daily update) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.963.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5069"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06374"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07779"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.288"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.647.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5449"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7830"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.013.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.440"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.961"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.048"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.898"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1171"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.884"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.237"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04981"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.261"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.185"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.541"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.359"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8944"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.580"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.133.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5449"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01555"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.555"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₈128"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.94%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.599"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8159"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.777.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4540"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05073"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  +0.00%  "
